$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Definite BOM" (sheet1) - only the selection changed in the source
# diff (no data/structural change), so just move the active selection.
# ---------------------------------------------------------------------------
$wsDefinite = $wb.Worksheets.Item("Definite BOM")
$wsDefinite.Range("B8").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Planned BOM" (sheet2) - replace the aluminum-rod / PVC-pipe "handle"
# parts with the newly found microcontroller + wire components.
# ---------------------------------------------------------------------------
$wsPlanned = $wb.Worksheets.Item("Planned BOM")

# Row 3: Sparkfun Qwiic Pro Micro (replaces the Leonardo/Xinput row)
$wsPlanned.Range("A3").Value = "Sparkfun Qwiic Pro Micro"

# Header row: add a "Link" column and rename the trailing "Qty" header to
# "Comments" (Qty now lives in column C).
$wsPlanned.Range("C1").Value = "Qty"
$wsPlanned.Range("D1").Value = "Link"
$wsPlanned.Range("E1").Value = "Comments"

$wsPlanned.Range("B3").Value = "Microcontroller"
$wsPlanned.Range("C3").Value = 1
$wsPlanned.Range("D3").Value = "https://www.microcenter.com/product/621833/sparkfun-electronics-qwiic-pro-micro-usb-c-(atmega32u4)"
$wsPlanned.Range("E3").Value = "Can be substituted for any Arduino Pro Micro or anything with an Atmega32U4 for HID support"

# Row 4: Stranded wire for the limit switch (replaces the Uno/Serial row)
$wsPlanned.Range("B4").Value = "Wire Limit Switch"
$wsPlanned.Range("C4").Value = 1
$wsPlanned.Range("D4").Value = "https://www.amazon.com/gp/product/B07T4SYVYG/"
$wsPlanned.Range("A4").Value = "Stranded Wire"
$wsPlanned.Range("E4").ClearContents()

# Old rows 5 & 6 (aluminum rod / PVC pipe parts) are gone entirely.
$wsPlanned.Range("A5:F6").ClearContents()

# Column D needs to be wide enough to show the long product link.
$wsPlanned.Columns.Item(4).ColumnWidth = 35.1666666

# Page setup was touched (orientation saved as portrait) when the sheet was
# last printed/previewed.
$wsPlanned.PageSetup.Orientation = 1

# Leave the cursor parked below the new table, matching the saved selection.
$wsPlanned.Range("C10").Select() | Out-Null
